$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old speedup formula that used to live in B6
$ws.Range("B6").ClearContents()

# --- Benchmark block: sample 1000 trajectories ---
$ws.Range("A6").Value = "16 actors, collect 100 trajectories, sample 1000 trajectories, deque of 5000 "
$ws.Range("A7").Value = "Ray"
$ws.Range("B7").Value = "Round 500, Training Time: 148.66"
$ws.Range("A8").Value = "LF"
$ws.Range("B8").Value = "Round 500, Training Time: 95.09"

# --- Benchmark block: sample 2000 trajectories ---
$ws.Range("A10").Value = "16 actors, collect 100 trajectories, sample 2000 trajectories, deque of 5000 "
$ws.Range("A11").Value = "Ray"
$ws.Range("A12").Value = "LF"
$ws.Range("B12").Value = "Round 500, Training Time: 95.97"
$ws.Range("B11").Value = "Round 500, Training Time: 191.55"

# Update the selected cell, mimicking the cursor position at save time
$ws.Range("B15").Select()
